{"js": "// 1. Version bump\nconst verResults = context.document.body.search(\"Version 2.2\", { matchCase: true });\nverResults.load(\"text\");\nawait context.sync();\nverResults.items[0].insertText(\"Version 2.3\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2. Data layer sentence rewrite\nconst dlResults = context.document.body.search(\n  \"The data layer also provides singletons that give all components access to settings and the logger.\",\n  { matchCase: true }\n);\ndlResults.load(\"text\");\nawait context.sync();\ndlResults.items[0].insertText(\n  \"The data layer also gives all components access to a logger and a settings singleton.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 3. Move the _GoBack bookmark to the end of the \"Version 2.3\" paragraph\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst verParaResults = context.document.body.search(\"Version 2.3\", { matchCase: true });\nverParaResults.load(\"text\");\nawait context.sync();\nconst endRange = verParaResults.items[0].getRange(Word.RangeLocation.end);\nendRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n\nreturn \"done\";\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Version bump\n$rng = $d.Content\n$rng.Find.Execute(\"Version 2.2\", $false, $false, $false, $false, $false, $true, 1, $false, \"Version 2.3\", 2)\n\n# 2. Data layer sentence rewrite\n$old = \"The data layer also provides singletons that give all components access to settings and the logger.\"\n$new = \"The data layer also gives all components access to a logger and a settings singleton.\"\n$rng2 = $d.Content\n$rng2.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n\n# 3. Move bookmark\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$bm.Delete()\n\n$txt = $d.Content.Text\n$idx = $txt.IndexOf(\"Version 2.3\")\n$endPos = $idx + 11  # length of \"Version 2.3\"\n$newRange = $d.Range($endPos, $endPos)\n$d.Bookmarks.Add(\"_GoBack\", $newRange)\n\nWrite-Output \"Exists: $($d.Bookmarks.Exists('_GoBack'))\"\n"}
